$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the header row (row 1), pushing existing
# data rows down by two.
$ws.Rows.Item(2).Resize(2).Insert()

# Populate the two newly inserted rows with the new accelerometer sample
# values.
$ws.Cells.Item(2, 1).Value = -0.9255759716033934
$ws.Cells.Item(2, 2).Value = 1.071177214384079
$ws.Cells.Item(2, 3).Value = -0.401322513818741

$ws.Cells.Item(3, 1).Value = -1.280259013175964
$ws.Cells.Item(3, 2).Value = 0.5940434336662286
$ws.Cells.Item(3, 3).Value = -1.134497284889223

# Drop the last three rows of the original data (now shifted to rows 22-24)
# so the sheet ends at row 21.
$ws.Rows.Item(22).Resize(3).Delete()
